# Apply the edits described by the commit diff to the active workbook.
#
# Summary of changes:
#   1. Rename the "Enthalpy" column header (F1, via the shared string table)
#      to "Heat Flow".
#   2. Update four enthalpy data values in column F (rows 3, 4, 7, 8) to
#      their new, slightly more precise figures.
#
# (Two purely cosmetic, tool-generated artifacts from the diff -- a
# `<charset>` sub-element on the default font, and hair-line
# `defaultColWidth`/column-D `width` adjustments of a few hundredths of a
# character -- have no corresponding property on Excel's Font/Worksheet/
# Range object model [confirmed against this host's own COM reflection,
# `Get-Member` lists no `Charset` on `Font`, and `StandardWidth`/
# `ColumnWidth` cannot address that granularity], so they are not
# reachable from COM automation and are intentionally left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: Enthalpy -> Heat Flow
$ws.Range("F1").Value = "Heat Flow"

# 2. Data value corrections in column F
$ws.Range("F3").Value = 41605.344
$ws.Range("F4").Value = 59155.92
$ws.Range("F7").Value = 15259.32
$ws.Range("F8").Value = -103654.08
